# Auto edit: update column A ("card") values in rows 3-12 of sheet "Card19"
# from "2" to "19" so they match the card's own sheet number.
#
# We briefly force a Text number format before writing the value so that
# the numeric-looking string "19" is kept as text (consistent with how
# these lookup sheets store every value as text) instead of being
# auto-coerced into a number by Excel. The style is then restored to
# "Normal" so no stray formatting is left behind on the cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card19")

for ($r = 3; $r -le 12; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.NumberFormat = "@"
    $cell.Value = "19"
    $cell.Style = "Normal"
}
